# Doing Updates for Financials
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACRX")

# --- Row 20: Total Other Income/Expenses Net ---
$ws.Range("D20").Value = -10200
$ws.Range("E20").Value = -8500
$ws.Range("F20").Value = -700

# --- Row 21: Earnings Before Interest And Taxes ---
$ws.Range("D21").Value = -47100
$ws.Range("E21").Value = -38400
$ws.Range("F21").Value = -18700
$ws.Range("J21").Value = "NA"

# --- Row 22: Interest Expense ---
$ws.Range("D22").Value = 3300
$ws.Range("E22").Value = 2800
$ws.Range("F22").Value = 3000

# --- Row 32: Other Items ---
$ws.Range("D32").Value = 10200
$ws.Range("E32").Value = 8500
$ws.Range("F32").Value = 700

# --- Cash Flow Statement: J column cells becoming "NA" ---
$ws.Range("J83").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"
